{"js": "// Applies the answer-key updates to the addition/subtraction practice\n// table: each cell in the 20x5 table holds a single arithmetic\n// expression (a single run of text) that gets replaced by its updated\n// version. Old values are all unique in the document, so a literal\n// search + replace for each pair is unambiguous.\nconst replacements = [\n  [\"50+5=55\", \"19+39=58\"],\n  [\"63+8=71\", \"84-2=82\"],\n  [\"55-37=18\", \"36-12=24\"],\n  [\"90-12=78\", \"92-18=74\"],\n  [\"78-24=54\", \"11+55=66\"],\n  [\"98-57=41\", \"33-28=5\"],\n  [\"17+1=18\", \"72-12=60\"],\n  [\"79-35=44\", \"88-17=71\"],\n  [\"73+16=89\", \"1+79=80\"],\n  [\"33+53=86\", \"27+3=30\"],\n  [\"79-28=51\", \"79+18=97\"],\n  [\"48-12=36\", \"82-66=16\"],\n  [\"47+10=57\", \"79-42=37\"],\n  [\"66-46=20\", \"23+76=99\"],\n  [\"97-76=21\", \"41-27=14\"],\n  [\"86+4=90\", \"9+76=85\"],\n  [\"26-5=21\", \"2+46=48\"],\n  [\"92-27=65\", \"87-6=81\"],\n  [\"41+3=44\", \"93-2=91\"],\n  [\"69+0=69\", \"61-58=3\"],\n  [\"96-74=22\", \"73+12=85\"],\n  [\"49+14=63\", \"42+43=85\"],\n  [\"74-68=6\", \"42+51=93\"],\n  [\"24-21=3\", \"87-22=65\"],\n  [\"60-19=41\", \"63-42=21\"],\n  [\"35+48=83\", \"65-35=30\"],\n  [\"36+17=53\", \"64+15=79\"],\n  [\"88-56=32\", \"61-44=17\"],\n  [\"90-7=83\", \"39+53=92\"],\n  [\"76+4=80\", \"80+1=81\"],\n  [\"65+34=99\", \"18+45=63\"],\n  [\"9+72=81\", \"93-38=55\"],\n  [\"52-20=32\", \"80+6=86\"],\n  [\"52+24=76\", \"96-40=56\"],\n  [\"56+34=90\", \"49-22=27\"],\n  [\"89-86=3\", \"80+16=96\"],\n  [\"47+27=74\", \"89-88=1\"],\n  [\"41-4=37\", \"33+34=67\"],\n  [\"72+3=75\", \"75-72=3\"],\n  [\"19+27=46\", \"63-2=61\"],\n  [\"58-25=33\", \"99-56=43\"],\n  [\"74-5=69\", \"74+25=99\"],\n  [\"72-65=7\", \"50-41=9\"],\n  [\"71-3=68\", \"94+5=99\"],\n  [\"43-36=7\", \"64-29=35\"],\n  [\"14-6=8\", \"25-10=15\"],\n  [\"68+24=92\", \"76-27=49\"],\n  [\"54-46=8\", \"6+66=72\"],\n  [\"22-0=22\", \"41+7=48\"],\n  [\"8+28=36\", \"69+15=84\"],\n  [\"83-20=63\", \"92-40=52\"],\n  [\"18+52=70\", \"90+4=94\"],\n  [\"74-25=49\", \"43-35=8\"],\n  [\"52+4=56\", \"44+34=78\"],\n  [\"82-57=25\", \"96-86=10\"],\n  [\"55-31=24\", \"80-18=62\"],\n  [\"26+50=76\", \"87-19=68\"],\n  [\"22+75=97\", \"35+6=41\"],\n  [\"34+17=51\", \"9+53=62\"],\n  [\"93-22=71\", \"23+69=92\"],\n  [\"18+79=97\", \"28+27=55\"],\n  [\"73-37=36\", \"87-58=29\"],\n  [\"88-33=55\", \"41-34=7\"],\n  [\"19+47=66\", \"50-23=27\"],\n  [\"2+36=38\", \"23-3=20\"],\n  [\"30+35=65\", \"42-6=36\"],\n  [\"85-0=85\", \"25+26=51\"],\n  [\"50+37=87\", \"28-4=24\"],\n  [\"76-13=63\", \"82-18=64\"],\n  [\"60-1=59\", \"49+12=61\"],\n  [\"67-32=35\", \"29-14=15\"],\n  [\"1+80=81\", \"72-28=44\"],\n  [\"57+8=65\", \"70-33=37\"],\n  [\"60+3=63\", \"37+50=87\"],\n  [\"63-7=56\", \"46+3=49\"],\n  [\"95+3=98\", \"37+58=95\"],\n  [\"55-48=7\", \"16+24=40\"],\n  [\"82+4=86\", \"74+24=98\"],\n  [\"18+53=71\", \"67+27=94\"],\n  [\"35-6=29\", \"6+47=53\"],\n  [\"54+8=62\", \"59-5=54\"],\n  [\"30+19=49\", \"86-34=52\"],\n  [\"29+8=37\", \"23-9=14\"],\n  [\"16+15=31\", \"94-71=23\"],\n  [\"53-25=28\", \"32-20=12\"],\n  [\"32-24=8\", \"78+11=89\"],\n  [\"65+11=76\", \"5+40=45\"],\n  [\"64-22=42\", \"21+58=79\"],\n  [\"92-34=58\", \"8-0=8\"],\n  [\"41+49=90\", \"51-25=26\"],\n  [\"50+20=70\", \"15+31=46\"],\n  [\"43-7=36\", \"45-20=25\"],\n  [\"90-78=12\", \"83-58=25\"],\n  [\"98-76=22\", \"6+76=82\"],\n  [\"92-6=86\", \"55+44=99\"],\n  [\"94-52=42\", \"54-38=16\"],\n  [\"96-12=84\", \"71+1=72\"],\n  [\"93+3=96\", \"63-33=30\"],\n  [\"6+65=71\", \"42-29=13\"],\n  [\"4-3=1\", \"94-64=30\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false, matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  // Every cell value in this table is unique, so exactly one match is expected;\n  // replace just the first (only) hit to avoid touching unrelated text.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Applies the answer-key updates to the addition/subtraction practice\n# table: each cell in the 20x5 table holds a single arithmetic\n# expression that gets replaced by its updated version via Find/Replace.\n# All old values are unique in the document, so each Find.Execute call\n# targets exactly one cell.\n\n$d = $word.ActiveDocument\n$wdReplaceAll = 2\n\nfunction ReplaceOne($oldText, $newText) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.MatchCase = $true\n  $find.MatchWildcards = $false\n  $find.Replacement.Text = $newText\n  $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $wdReplaceAll) | Out-Null\n}\n\nReplaceOne \"50+5=55\" \"19+39=58\"\nReplaceOne \"63+8=71\" \"84-2=82\"\nReplaceOne \"55-37=18\" \"36-12=24\"\nReplaceOne \"90-12=78\" \"92-18=74\"\nReplaceOne \"78-24=54\" \"11+55=66\"\nReplaceOne \"98-57=41\" \"33-28=5\"\nReplaceOne \"17+1=18\" \"72-12=60\"\nReplaceOne \"79-35=44\" \"88-17=71\"\nReplaceOne \"73+16=89\" \"1+79=80\"\nReplaceOne \"33+53=86\" \"27+3=30\"\nReplaceOne \"79-28=51\" \"79+18=97\"\nReplaceOne \"48-12=36\" \"82-66=16\"\nReplaceOne \"47+10=57\" \"79-42=37\"\nReplaceOne \"66-46=20\" \"23+76=99\"\nReplaceOne \"97-76=21\" \"41-27=14\"\nReplaceOne \"86+4=90\" \"9+76=85\"\nReplaceOne \"26-5=21\" \"2+46=48\"\nReplaceOne \"92-27=65\" \"87-6=81\"\nReplaceOne \"41+3=44\" \"93-2=91\"\nReplaceOne \"69+0=69\" \"61-58=3\"\nReplaceOne \"96-74=22\" \"73+12=85\"\nReplaceOne \"49+14=63\" \"42+43=85\"\nReplaceOne \"74-68=6\" \"42+51=93\"\nReplaceOne \"24-21=3\" \"87-22=65\"\nReplaceOne \"60-19=41\" \"63-42=21\"\nReplaceOne \"35+48=83\" \"65-35=30\"\nReplaceOne \"36+17=53\" \"64+15=79\"\nReplaceOne \"88-56=32\" \"61-44=17\"\nReplaceOne \"90-7=83\" \"39+53=92\"\nReplaceOne \"76+4=80\" \"80+1=81\"\nReplaceOne \"65+34=99\" \"18+45=63\"\nReplaceOne \"9+72=81\" \"93-38=55\"\nReplaceOne \"52-20=32\" \"80+6=86\"\nReplaceOne \"52+24=76\" \"96-40=56\"\nReplaceOne \"56+34=90\" \"49-22=27\"\nReplaceOne \"89-86=3\" \"80+16=96\"\nReplaceOne \"47+27=74\" \"89-88=1\"\nReplaceOne \"41-4=37\" \"33+34=67\"\nReplaceOne \"72+3=75\" \"75-72=3\"\nReplaceOne \"19+27=46\" \"63-2=61\"\nReplaceOne \"58-25=33\" \"99-56=43\"\nReplaceOne \"74-5=69\" \"74+25=99\"\nReplaceOne \"72-65=7\" \"50-41=9\"\nReplaceOne \"71-3=68\" \"94+5=99\"\nReplaceOne \"43-36=7\" \"64-29=35\"\nReplaceOne \"14-6=8\" \"25-10=15\"\nReplaceOne \"68+24=92\" \"76-27=49\"\nReplaceOne \"54-46=8\" \"6+66=72\"\nReplaceOne \"22-0=22\" \"41+7=48\"\nReplaceOne \"8+28=36\" \"69+15=84\"\nReplaceOne \"83-20=63\" \"92-40=52\"\nReplaceOne \"18+52=70\" \"90+4=94\"\nReplaceOne \"74-25=49\" \"43-35=8\"\nReplaceOne \"52+4=56\" \"44+34=78\"\nReplaceOne \"82-57=25\" \"96-86=10\"\nReplaceOne \"55-31=24\" \"80-18=62\"\nReplaceOne \"26+50=76\" \"87-19=68\"\nReplaceOne \"22+75=97\" \"35+6=41\"\nReplaceOne \"34+17=51\" \"9+53=62\"\nReplaceOne \"93-22=71\" \"23+69=92\"\nReplaceOne \"18+79=97\" \"28+27=55\"\nReplaceOne \"73-37=36\" \"87-58=29\"\nReplaceOne \"88-33=55\" \"41-34=7\"\nReplaceOne \"19+47=66\" \"50-23=27\"\nReplaceOne \"2+36=38\" \"23-3=20\"\nReplaceOne \"30+35=65\" \"42-6=36\"\nReplaceOne \"85-0=85\" \"25+26=51\"\nReplaceOne \"50+37=87\" \"28-4=24\"\nReplaceOne \"76-13=63\" \"82-18=64\"\nReplaceOne \"60-1=59\" \"49+12=61\"\nReplaceOne \"67-32=35\" \"29-14=15\"\nReplaceOne \"1+80=81\" \"72-28=44\"\nReplaceOne \"57+8=65\" \"70-33=37\"\nReplaceOne \"60+3=63\" \"37+50=87\"\nReplaceOne \"63-7=56\" \"46+3=49\"\nReplaceOne \"95+3=98\" \"37+58=95\"\nReplaceOne \"55-48=7\" \"16+24=40\"\nReplaceOne \"82+4=86\" \"74+24=98\"\nReplaceOne \"18+53=71\" \"67+27=94\"\nReplaceOne \"35-6=29\" \"6+47=53\"\nReplaceOne \"54+8=62\" \"59-5=54\"\nReplaceOne \"30+19=49\" \"86-34=52\"\nReplaceOne \"29+8=37\" \"23-9=14\"\nReplaceOne \"16+15=31\" \"94-71=23\"\nReplaceOne \"53-25=28\" \"32-20=12\"\nReplaceOne \"32-24=8\" \"78+11=89\"\nReplaceOne \"65+11=76\" \"5+40=45\"\nReplaceOne \"64-22=42\" \"21+58=79\"\nReplaceOne \"92-34=58\" \"8-0=8\"\nReplaceOne \"41+49=90\" \"51-25=26\"\nReplaceOne \"50+20=70\" \"15+31=46\"\nReplaceOne \"43-7=36\" \"45-20=25\"\nReplaceOne \"90-78=12\" \"83-58=25\"\nReplaceOne \"98-76=22\" \"6+76=82\"\nReplaceOne \"92-6=86\" \"55+44=99\"\nReplaceOne \"94-52=42\" \"54-38=16\"\nReplaceOne \"96-12=84\" \"71+1=72\"\nReplaceOne \"93+3=96\" \"63-33=30\"\nReplaceOne \"6+65=71\" \"42-29=13\"\nReplaceOne \"4-3=1\" \"94-64=30\"\n"}
